$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new technology rows to the table
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "hs_house"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "dh_grid"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "po_turbine"

# Grow Table1 so it covers the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B8"))

# Update the active selection to reflect where the user left off editing
$ws.Range("B9").Select()
